$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -ne "Picture 4") {
        $shape.Delete()
    }
}
